# Update test xlsx files to use ID (NIT and DUI)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: "nit" -> "id"
$ws.Range("A1").Value = "id"

# Fill in the previously-blank row 4 with a new employee record
$ws.Range("A4").Value = 55544433
$ws.Range("B4").Value = "IJK LMN"
$ws.Range("C4").Value = 456
$ws.Range("D4").Value = "ijk"

# Give A4 its own number format + font (Calibri 11, black, left aligned)
$ws.Range("A4").NumberFormat = "00000000000000"
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").Font.Name = "Calibri"
$ws.Range("A4").Font.Size = 11
$ws.Range("A4").Font.Color = 0

# Row 4 is a touch taller than the rest
$ws.Range("A4").RowHeight = 13.8

# Move the active selection to A5
[void]$ws.Range("A5").Select()

Write-Host "done"
